$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F; this shifts F.. rightwards and
# auto-updates formula references.
$ws.Columns("F:F").Insert()

# New header for inserted column
$ws.Range("F1").Value = "incometax_tu"

# New formula for the inserted column (rows 2-7)
$ws.Range("F2").Formula = "=MIN(G2,H2)"
$ws.Range("F3:F7").Formula = "=MIN(G3,H3)"

# Re-group the formulas that were shifted out of F into G..P so the
# (pre-existing) shared-formula blocks survive the column insert.
$ws.Range("P3:P5").Formula = "=G3+I3"
$ws.Range("P7").Formula = "=G7+I7"
$ws.Range("M6:M7").Formula = "=N6"

# Update selection to match the recorded interactive state
$ws.Range("F12").Select()
